# Singapore Premier League workbook update
# The upstream data refresh re-paired several fixtures that share the same
# match date: for each of the row-pairs below, the two rows' match details
# (id, teams, score, result, odds, etc. -- i.e. everything except the
# leading row index in column A) were swapped with one another.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row1, row2) pairs whose B:AC contents (columns B and F..AC; C/D/E are
# identical within each pair so including them in the swap is a no-op)
# need to be exchanged.
$rowPairs = @(
    @(8, 9),
    @(10, 11),
    @(18, 19),
    @(20, 21),
    @(22, 23),
    @(28, 29),
    @(38, 39),
    @(43, 44),
    @(47, 48),
    @(51, 52)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
